# Scheduled market-price refresh for the Gilgamesh Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H, I, J, K, L, M, N)
# on a handful of rows across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets with freshly
# polled marketboard data. CUL is untouched this run.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Sheets.Item("ALC")
# row 3
$ws.Range("H3").Value = 18995.5
$ws.Range("J3").Value = 18995.5
$ws.Range("L3").Value = 18995.5
$ws.Range("N3").Value = -19223.5
# row 98
$ws.Range("H98").Value = 6378
$ws.Range("I98").Value = 6378
$ws.Range("K98").Value = 6378
$ws.Range("M98").Value = -4880
# row 99
$ws.Range("H99").Value = 313.8
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null
# row 102
$ws.Range("H102").Value = 18995.5
$ws.Range("J102").Value = 18995.5
$ws.Range("L102").Value = 18995.5
$ws.Range("N102").Value = -25485.5
# row 122
$ws.Range("H122").Value = 6378
$ws.Range("I122").Value = 6378
$ws.Range("K122").Value = 19134
$ws.Range("M122").Value = -16684

# --- ARM ---
$ws = $wb.Sheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 831.625
$ws.Range("J2").Value = 1268.2222
$ws.Range("L2").Value = 1268.2222
$ws.Range("N2").Value = -1494.2222
# row 61
$ws.Range("H61").Value = 3968.5652
$ws.Range("I61").Value = 2516.5
$ws.Range("J61").Value = 7287.5713
$ws.Range("K61").Value = 2516.5
$ws.Range("L61").Value = 7287.5713
$ws.Range("M61").Value = -2304.5
$ws.Range("N61").Value = -7711.5713
# row 74
$ws.Range("H74").Value = 266647.2
$ws.Range("I74").Value = 371867.53
$ws.Range("K74").Value = 371867.53
$ws.Range("M74").Value = -370993.53
# row 77
$ws.Range("H77").Value = 266647.2
$ws.Range("I77").Value = 371867.53
$ws.Range("K77").Value = 1859337.65
$ws.Range("M77").Value = -1854969.65
# row 116
$ws.Range("H116").Value = 831.625
$ws.Range("J116").Value = 1268.2222
$ws.Range("L116").Value = 1268.2222
$ws.Range("N116").Value = -5856.2222
# row 132
$ws.Range("H132").Value = 2651.9678
$ws.Range("J132").Value = 5449.9
$ws.Range("L132").Value = 16349.7
$ws.Range("N132").Value = -21409.7
# row 136
$ws.Range("H136").Value = 3968.5652
$ws.Range("I136").Value = 2516.5
$ws.Range("J136").Value = 7287.5713
$ws.Range("K136").Value = 7549.5
$ws.Range("L136").Value = 21862.7139
$ws.Range("M136").Value = -4999.5
$ws.Range("N136").Value = -26962.7139

# --- BSM ---
$ws = $wb.Sheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 831.625
$ws.Range("J3").Value = 1268.2222
$ws.Range("L3").Value = 1268.2222
$ws.Range("N3").Value = -1496.2222
# row 80
$ws.Range("H80").Value = 307.92307
$ws.Range("I80").Value = 108.2
$ws.Range("J80").Value = 432.75
$ws.Range("K80").Value = 108.2
$ws.Range("L80").Value = 432.75
$ws.Range("M80").Value = 889.8
$ws.Range("N80").Value = -2428.75
# row 83
$ws.Range("H83").Value = 307.92307
$ws.Range("I83").Value = 108.2
$ws.Range("J83").Value = 432.75
$ws.Range("K83").Value = 541
$ws.Range("L83").Value = 2163.75
$ws.Range("M83").Value = 4451
$ws.Range("N83").Value = -12147.75

# --- CRP ---
$ws = $wb.Sheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3537.537
$ws.Range("I31").Value = 2299.1025
$ws.Range("K31").Value = 2299.1025
$ws.Range("M31").Value = -2004.1025
# row 34
$ws.Range("H34").Value = 3537.537
$ws.Range("I34").Value = 2299.1025
$ws.Range("K34").Value = 2299.1025
$ws.Range("M34").Value = -2097.1025
# row 99
$ws.Range("H99").Value = 6133
$ws.Range("I99").Value = 4666.3335
$ws.Range("J99").Value = 7599.6665
$ws.Range("K99").Value = 4666.3335
$ws.Range("L99").Value = 7599.6665
$ws.Range("M99").Value = -3168.3335
$ws.Range("N99").Value = -10595.6665
# row 122
$ws.Range("H122").Value = 3290.9
$ws.Range("I122").Value = 2720.1428
$ws.Range("K122").Value = 8160.428400000001
$ws.Range("M122").Value = -5710.428400000001
# row 126
$ws.Range("H126").Value = 6133
$ws.Range("I126").Value = 4666.3335
$ws.Range("J126").Value = 7599.6665
$ws.Range("K126").Value = 13999.0005
$ws.Range("L126").Value = 22798.9995
$ws.Range("M126").Value = -11529.0005
$ws.Range("N126").Value = -27738.9995

# --- GSM ---
$ws = $wb.Sheets.Item("GSM")
# row 46
$ws.Range("H46").Value = 22460.334
$ws.Range("J46").Value = 33190.5
$ws.Range("L46").Value = 33190.5
$ws.Range("N46").Value = -33502.5
# row 70
$ws.Range("H70").Value = 162131.39
$ws.Range("I70").Value = 259912.62
$ws.Range("K70").Value = 259912.62
$ws.Range("M70").Value = -259642.62
# row 73
$ws.Range("H73").Value = 162131.39
$ws.Range("I73").Value = 259912.62
$ws.Range("K73").Value = 259912.62
$ws.Range("M73").Value = -258976.62
# row 102
$ws.Range("H102").Value = 1252
$ws.Range("I102").Value = 845
$ws.Range("K102").Value = 845
$ws.Range("M102").Value = 777
# row 122
$ws.Range("H122").Value = 5430.467
$ws.Range("I122").Value = 3973.3635
$ws.Range("J122").Value = 9437.5
$ws.Range("K122").Value = 11920.0905
$ws.Range("L122").Value = 28312.5
$ws.Range("M122").Value = -9470.0905
$ws.Range("N122").Value = -33212.5
# row 126
$ws.Range("H126").Value = 12156.833
$ws.Range("I126").Value = 7485.75
$ws.Range("J126").Value = 21499
$ws.Range("K126").Value = 22457.25
$ws.Range("L126").Value = 64497
$ws.Range("M126").Value = -19987.25
$ws.Range("N126").Value = -69437

# --- LTW ---
$ws = $wb.Sheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 739.8
$ws.Range("J22").Value = 1250
$ws.Range("L22").Value = 1250
$ws.Range("N22").Value = -1840
# row 27
$ws.Range("H27").Value = 739.8
$ws.Range("J27").Value = 1250
$ws.Range("L27").Value = 1250
$ws.Range("N27").Value = -1464
# row 40
$ws.Range("H40").Value = 256496
$ws.Range("I40").Value = 502999.5
$ws.Range("J40").Value = 9992.5
$ws.Range("K40").Value = 502999.5
$ws.Range("L40").Value = 9992.5
$ws.Range("M40").Value = -502863.5
$ws.Range("N40").Value = -10264.5
# row 122
$ws.Range("H122").Value = 3873.7778
$ws.Range("J122").Value = 3812.5
$ws.Range("L122").Value = 11437.5
$ws.Range("N122").Value = -16337.5
# row 132
$ws.Range("H132").Value = 25312.143
$ws.Range("I132").Value = 4590
$ws.Range("J132").Value = 33601
$ws.Range("K132").Value = 13770
$ws.Range("L132").Value = 100803
$ws.Range("M132").Value = -11240
$ws.Range("N132").Value = -105863

# --- WVR ---
$ws = $wb.Sheets.Item("WVR")
# row 96
$ws.Range("H96").Value = 3537.6924
$ws.Range("I96").Value = 3643
$ws.Range("J96").Value = 3414.8333
$ws.Range("K96").Value = 3643
$ws.Range("L96").Value = 3414.8333
$ws.Range("M96").Value = -2270
$ws.Range("N96").Value = -6160.8333
# row 107
$ws.Range("H107").Value = 687.6667
$ws.Range("I107").Value = 461
$ws.Range("J107").Value = 1005
$ws.Range("K107").Value = 1383
$ws.Range("L107").Value = 3015
$ws.Range("M107").Value = 537
$ws.Range("N107").Value = -6855
# row 119
$ws.Range("H119").Value = 77286.5
$ws.Range("J119").Value = 77286.5
$ws.Range("L119").Value = 77286.5
$ws.Range("N119").Value = -86962.5
# row 122
$ws.Range("H122").Value = 31251836
$ws.Range("I122").Value = 1874.75
$ws.Range("K122").Value = 5624.25
$ws.Range("M122").Value = -3174.25
# row 126
$ws.Range("H126").Value = 1193
$ws.Range("I126").Value = 1193
$ws.Range("K126").Value = 3579
$ws.Range("M126").Value = -1109
